$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1) Fix the typo "spint" -> "spin" in the OUTCAR argument example.
# ---------------------------------------------------------------
$d.Content.Find.Execute(
    "use spint=0, 1, 2, 3 (respectively). For example:",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "use spin=0, 1, 2, 3 (respectively). For example:", 2)

# ---------------------------------------------------------------
# 2) Register the additional ListLabel character styles (20-28)
#    that the new numbering definitions reference.
# ---------------------------------------------------------------
# wdStyleTypeCharacter = 2
function Add-ListLabelStyle($name, $latinFont, $csFont, $bold, $size) {
    $s = $d.Styles.Add($name, 2)
    $s.NameLocal = $name
    $s.QuickStyle = $true
    $f = $s.Font
    if ($latinFont) {
        $f.NameAscii = $latinFont
        $f.Name = $latinFont
    }
    if ($csFont) { $f.NameBi = $csFont }
    if ($bold) { $f.Bold = $true }
    if ($size) { $f.Size = $size }
}

Add-ListLabelStyle "ListLabel 20" "Times New Roman" "Wingdings" $true 11
Add-ListLabelStyle "ListLabel 21" $null "Courier New" $false $null
Add-ListLabelStyle "ListLabel 22" $null "Wingdings" $false $null
Add-ListLabelStyle "ListLabel 23" $null "Symbol" $false $null
Add-ListLabelStyle "ListLabel 24" $null "Courier New" $false $null
Add-ListLabelStyle "ListLabel 25" $null "Wingdings" $false $null
Add-ListLabelStyle "ListLabel 26" $null "Symbol" $false $null
Add-ListLabelStyle "ListLabel 27" $null "Courier New" $false $null
Add-ListLabelStyle "ListLabel 28" $null "Wingdings" $false $null
